$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @('D5','D6','D8','D9','D15','D18','D19','D20','D21','D24','D25','D27','D30','D31','D32','D34','D35','D36','D38','D39','D40','D41','D44','D45','D47','D48','D49','D50','D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '59.616.33'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '2.614.14'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('D5').Value = '537.40'
$ws.Range('E5').Value = '  +2.44%  '
$ws.Range('D6').Value = '142.20'
$ws.Range('E6').Value = '  +1.61%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').Value = '0.566'
$ws.Range('E8').Value = '  +0.32%  '
$ws.Range('D9').Value = '6.58'
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('E10').Value = '  +0.86%  '
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('D13').Value = '3.079.49'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').Value = '59.548.61'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').Value = '20.74'
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('D16').Value = '2.634.06'
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').Value = '340.35'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('D19').Value = '4.35'
$ws.Range('E19').Value = '  +0.78%  '
$ws.Range('D20').Value = '10.11'
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').Value = '6.34'
$ws.Range('E21').Value = '  -1.61%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('B24').Value = 'Kaspa'
$ws.Range('C24').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D24').Value = '0.166'
$ws.Range('E24').Value = '  -1.55%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').Value = '0.408'
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('D27').Value = '7.22'
$ws.Range('E27').Value = '  +1.93%  '
$ws.Range('D28').Value = '0.0₃0744'
$ws.Range('E28').Value = '  +2.43%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '1.66'
$ws.Range('E30').Value = '  +4.60%  '
$ws.Range('D31').Value = '5.83'
$ws.Range('E31').Value = '  -1.99%  '
$ws.Range('D32').Value = '18.79'
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('E33').Value = '  +0.94%  '
$ws.Range('D34').Value = '3.98'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').Value = '1.12'
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('D36').Value = '0.833'
$ws.Range('E36').Value = '  +2.51%  '
$ws.Range('E37').Value = '  -1.95%  '
$ws.Range('D38').Value = '0.826'
$ws.Range('E38').Value = '  -0.69%  '
$ws.Range('D39').Value = '3.53'
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').Value = '277.89'
$ws.Range('E40').Value = '  +1.84%  '
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('D44').Value = '0.0948'
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('D45').Value = '0.0525'
$ws.Range('E45').Value = '  +1.94%  '
$ws.Range('D46').Value = '1.954.33'
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('D47').Value = '0.0223'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '18.40'
$ws.Range('E48').Value = '  +0.40%  '
$ws.Range('D49').Value = '4.51'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').Value = '111.69'
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('D51').Value = '4.75'
$ws.Range('E51').Value = '  +0.69%  '
